# Refresh the cryptos list (prices + 1h volume deltas), matching the
# automated "Updated cryptos list ... with GitHub Actions" commit.
# Price cells in column D are text (e.g. "59.691.43", "0.620") so a
# leading apostrophe is used to force Excel to keep them as text instead
# of silently coercing them to numbers (which would mangle values like
# "2.644.85" or drop significant trailing zeros like "0.620" -> 0.62).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.691.43"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "'2.648.23"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'518.66"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").Value = "'147.06"
$ws.Range("E6").Value = "  -1.14%  "

$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("E8").Value = "  +0.67%  "

$ws.Range("D9").Value = "'2.672.58"
$ws.Range("E9").Value = "  +2.20%  "

$ws.Range("E10").Value = "  +3.10%  "

$ws.Range("E11").Value = "  +2.13%  "

$ws.Range("D12").Value = "'0.341"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("E13").Value = "  -1.43%  "

$ws.Range("D14").Value = "'3.107.90"
$ws.Range("E14").Value = "  +1.32%  "

$ws.Range("D15").Value = "'59.540.55"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").Value = "'21.32"
$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").Value = "'2.660.89"
$ws.Range("E18").Value = "  +1.83%  "

$ws.Range("D19").Value = "'4.63"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").Value = "'346.58"
$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("D21").Value = "'10.53"
$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("D24").Value = "'61.08"
$ws.Range("E24").Value = "  +1.48%  "

$ws.Range("D25").Value = "'0.425"
$ws.Range("E25").Value = "  +1.55%  "

$ws.Range("D26").Value = "'2.764.85"
$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").Value = "'0.162"
$ws.Range("E28").Value = "  -0.38%  "

$ws.Range("D29").Value = "'0.0₃0825"
$ws.Range("E29").Value = "  +2.60%  "

$ws.Range("E30").Value = "  +2.39%  "

$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("E32").Value = "  +8.39%  "

$ws.Range("D33").Value = "'19.04"
$ws.Range("E33").Value = "  +0.68%  "

$ws.Range("E34").Value = "  -0.24%  "

$ws.Range("D35").Value = "'1.08"
$ws.Range("E35").Value = "  +19.22%  "

$ws.Range("D36").Value = "'149.64"
$ws.Range("E36").Value = "  -0.32%  "

$ws.Range("D37").Value = "'4.05"
$ws.Range("E37").Value = "  +2.54%  "

$ws.Range("E38").Value = "  +2.50%  "

$ws.Range("D39").Value = "'0.873"
$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("D40").Value = "'36.75"
$ws.Range("E40").Value = "  +1.11%  "

$ws.Range("E41").Value = "  +3.64%  "

$ws.Range("D42").Value = "'1.44"
$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").Value = "'285.23"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").Value = "'0.620"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").Value = "'0.0997"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'19.93"
$ws.Range("E46").Value = "  +2.33%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'0.992"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("E48").Value = "  +0.25%  "

$ws.Range("D49").Value = "'0.0234"
$ws.Range("E49").Value = "  +0.64%  "

$ws.Range("D50").Value = "'4.75"
$ws.Range("E50").Value = "  +1.61%  "

$ws.Range("E51").Value = "  -1.25%  "
